# Edit script for RE-L00-Organization.pptx
# Updates the "Lectures" schedule (slide 9) and "Exercises" schedule (slide 10)
# to reflect the altered timetable (21.12.2022 sessions cancelled, later
# sessions shifted back by one week, and a new final Exercises entry added).
#
# NOTE: this COM host's expression parser misbehaves when a parenthesized
# expression directly follows another call/literal used as a positional
# argument (e.g. "Foo $x.Bar(1,2) (\"str\" + $y)" silently drops args, and
# "Foo 13 (\"a\" + \"b\")" throws "Cannot invoke ... System.Int64"). To stay
# safe, every intermediate value (paragraph objects, concatenated strings)
# is always assigned to a variable first and only bare variables are passed
# into functions/methods.

$p = $ppt.ActivePresentation

# Useful special characters (built from code points to dodge any encoding
# pitfalls going through the COM-interop bridge).
$arrow = [char]0x2192   # "->" ... actually U+2192 RIGHTWARDS ARROW "→"
$ndash = [char]0x2013   # U+2013 EN DASH "–"

function Set-WholeRunText {
    # Replace the full text of a paragraph with brand-new text as a single
    # run, while preserving the run's existing formatting (rPr). Directly
    # assigning the final text can make the host diff old vs. new text and
    # keep a stale run boundary (or, if we go through an empty string, drop
    # the rPr entirely) -- routing through an unrelated placeholder string
    # first sidesteps both problems because it shares no common prefix or
    # suffix with either the old or the new text.
    param($para, $newText)

    $para.Text = "QQPLACEHOLDERQQ"
    $para.Text = $newText
}

function Set-SuffixRunText {
    # Split "<unchanged prefix>" + "<old suffix>" into two runs:
    # "<unchanged prefix>" (left alone) and a new run containing $newSuffix,
    # by overwriting only the characters after the prefix.
    param($para, $prefixLen, $newSuffix)

    $oldLen = $para.Length
    $subLen = $oldLen - $prefixLen
    $startPos = $prefixLen + 1
    $sub = $para.Characters($startPos, $subLen)
    $sub.Text = $newSuffix
}

# ---------------------------------------------------------------------------
# Slide 9 ("Lectures" schedule)
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(2)
$tr9 = $sh9.TextFrame.TextRange

$p9_9 = $tr9.Paragraphs(9, 1)
$text9_9 = "21.12.2022 " + $arrow + " No Lecture"
Set-WholeRunText $p9_9 $text9_9

$p9_10 = $tr9.Paragraphs(10, 1)
$text9_10 = "Documentation " + $ndash + " Formal Requirements Specification (L09)"
Set-SuffixRunText $p9_10 13 $text9_10

$p9_11 = $tr9.Paragraphs(11, 1)
$text9_11 = "Validation (L10)"
Set-SuffixRunText $p9_11 13 $text9_11

$p9_12 = $tr9.Paragraphs(12, 1)
$text9_12 = "Traceability (L11)"
Set-SuffixRunText $p9_12 13 $text9_12

$p9_13 = $tr9.Paragraphs(13, 1)
$text9_13 = "Requirements Management (L12) and Tool Support (L13)"
Set-SuffixRunText $p9_13 13 $text9_13

# ---------------------------------------------------------------------------
# Slide 10 ("Exercises" schedule)
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange

$p10_5 = $tr10.Paragraphs(5, 1)
$text10_5 = "21.12.2022 " + $arrow + " No Exercise"
Set-WholeRunText $p10_5 $text10_5

$p10_6 = $tr10.Paragraphs(6, 1)
$text10_6 = "Exercise 05 " + $ndash + " Coloured Petri Nets I"
Set-SuffixRunText $p10_6 13 $text10_6

$p10_7 = $tr10.Paragraphs(7, 1)
$text10_7 = "Exercise 06 " + $ndash + " Coloured Petri Nets II"
Set-SuffixRunText $p10_7 13 $text10_7

$p10_8 = $tr10.Paragraphs(8, 1)
$text10_8 = "Bonus Task"
Set-SuffixRunText $p10_8 13 $text10_8

# New trailing entry, shifted down from the old "25.01.2023" line.
$newParaText = "`r01.02.2023 " + $arrow + " Exercise 07 " + $ndash + " Management and Traceability (MC) "
$null = $tr10.InsertAfter($newParaText)
